# "Generate Report for Handoff"
# Flips the localization-status report from "In Translation" to
# "Ready for handoff" and refreshes the associated timestamps, widening
# the status/date columns so the new (longer) text fits.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"          # zh-cn status
$ws.Range("F2").Value = "Ready for handoff"          # de-de status
$ws.Range("G2").Value = "2016-08-28 00:56:04"        # Latest HO Xliff Generate Date
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss" # keep original datetime display format
$ws.Columns.Item(5).ColumnWidth = 16.25               # widen zh-cn status col
$ws.Columns.Item(6).ColumnWidth = 16.25               # widen de-de status col

# --- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"          # Status
$ws.Range("H2").Value = "2016-08-28 00:55:57"        # Latest Handoff Datetime
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss" # keep original datetime display format
$ws.Columns.Item(3).ColumnWidth = 16.25               # widen Status col

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"          # Status
$ws.Range("H2").Value = "2016-08-28 00:56:04"        # Latest Handoff Datetime
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss" # keep original datetime display format
$ws.Columns.Item(3).ColumnWidth = 16.25               # widen Status col
